# Apply odds updates to "Jogos da Semana" worksheet (row 5 and row 11)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 5 updates ---
$ws.Range("I5").Value = 6
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("Q5").Value = 2.4
$ws.Range("R5").Value = 1.53
$ws.Range("S5").Value = 1.53
$ws.Range("T5").Value = 2.38
$ws.Range("U5").Value = 2.5
$ws.Range("V5").Value = 1.5
$ws.Range("W5").Value = 5
$ws.Range("AC5").Value = 7
$ws.Range("AF5").Value = 101
$ws.Range("AJ5").Value = 67
$ws.Range("AO5").Value = 8.5
$ws.Range("AP5").Value = 26
$ws.Range("AS5").Value = 251
$ws.Range("AT5").Value = 2.38
$ws.Range("AU5").Value = 10
$ws.Range("AV5").Value = 81
$ws.Range("AY5").Value = 41
$ws.Range("AZ5").Value = 151

# --- Row 11 updates ---
$ws.Range("H11").Value = 3.6
$ws.Range("I11").Value = 1.83
$ws.Range("K11").Value = 2.4
$ws.Range("L11").Value = 2.22
$ws.Range("AC11").Value = 12
$ws.Range("AD11").Value = 7.1
$ws.Range("AG11").Value = 8.25
$ws.Range("AJ11").Value = 15.5
$ws.Range("AK11").Value = 13.5
$ws.Range("AM11").Value = 400
$ws.Range("AU11").Value = 6.5
$ws.Range("AX11").Value = 8.25
$ws.Range("AZ11").Value = 25

$wb.Save()
